# Update cryptocurrency price/volume symbol list (refreshed data snapshot).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Values in D and E are
# stored as literal text in the workbook (e.g. "305.58", "1.80%"), so a
# leading apostrophe is used to force Excel to keep them as text instead
# of auto-converting to numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''305.58'
$ws.Range('E2').Value = '''1.80%'
$ws.Range('D3').Value = '''35.83'
$ws.Range('E3').Value = '''1.46%'
$ws.Range('D4').Value = '''5.066'
$ws.Range('E4').Value = '''0.33%'
$ws.Range('D5').Value = '''0.08056'
$ws.Range('E5').Value = '''1.40%'
$ws.Range('D6').Value = '''1.923'
$ws.Range('E6').Value = '''1.49%'
$ws.Range('D7').Value = '''4.166'
$ws.Range('E7').Value = '''3.35%'
$ws.Range('E8').Value = '''0.89%'
$ws.Range('D9').Value = '''0.9289'
$ws.Range('E9').Value = '''0.18%'
$ws.Range('D10').Value = '''0.1336'
$ws.Range('E10').Value = '''-1.08%'
$ws.Range('D11').Value = '''0.1903'
$ws.Range('E11').Value = '''0.22%'
$ws.Range('D12').Value = '''0.09191'
$ws.Range('E12').Value = '''1.08%'
$ws.Range('D13').Value = '''0.03472'
$ws.Range('E13').Value = '''1.17%'
$ws.Range('D14').Value = '''0.09909'
$ws.Range('E14').Value = '''-0.08%'
$ws.Range('D15').Value = '''0.001416'
$ws.Range('E15').Value = '''1.71%'
$ws.Range('D16').Value = '''0.006701'
$ws.Range('E16').Value = '''14.45%'
$ws.Range('E17').Value = '''2.38%'
$ws.Range('E18').Value = '''1.74%'
$ws.Range('E19').Value = '''0.58%'
$ws.Range('E20').Value = '''3.48%'
$ws.Range('D21').Value = '''5.177'
$ws.Range('E21').Value = '''2.46%'
$ws.Range('D22').Value = '''0.2538'
$ws.Range('E22').Value = '''5.87%'
$ws.Range('D23').Value = '''0.04420'
$ws.Range('E23').Value = '''-1.60%'
$ws.Range('D24').Value = '''0.001236'
$ws.Range('E24').Value = '''1.89%'
$ws.Range('D25').Value = '''0.004699'
$ws.Range('E25').Value = '''-1.38%'
$ws.Range('D26').Value = '''0.0001300'
$ws.Range('E26').Value = '''5.62%'
$ws.Range('D27').Value = '''0.0003137'
$ws.Range('E27').Value = '''4.49%'
$ws.Range('D39').Value = '''0.01989'
$ws.Range('E39').Value = '''5.25%'
$ws.Range('D40').Value = '''0.05145'
$ws.Range('E40').Value = '''8.48%'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Value = '''0.007625'
$ws.Range('E41').Value = '''3.93%'
$ws.Range('B42').Value = 'Dexo'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range('D42').Value = '''0.01034'
$ws.Range('E42').Value = '''0.44%'
$ws.Range('D43').Value = '''0.1365'
$ws.Range('E43').Value = '''3.10%'
$ws.Range('D44').Value = '''0.002100'
$ws.Range('E44').Value = '''-0.55%'
$ws.Range('D45').Value = '''0.01074'
$ws.Range('E45').Value = '''-2.44%'
$ws.Range('D46').Value = '''0.00006321'
$ws.Range('E46').Value = '''0.60%'
$ws.Range('D47').Value = '''0.00000000750'
$ws.Range('E47').Value = '''-0.04%'
$ws.Range('E48').Value = '''-1.70%'
$ws.Range('E49').Value = '''-3.32%'
$ws.Range('D50').Value = '''0.00002100'
$ws.Range('E50').Value = '''-0.04%'
$ws.Range('D51').Value = '''0.0002000'
$ws.Range('E51').Value = '''-0.04%'
